$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Goals Completed:" column (D) added to the Week2 sheet.
$ws.Columns.Item(4).ColumnWidth = 20.140625

# Rows 4 and 6 grow taller to fit the newly entered text.
$ws.Rows.Item(4).RowHeight = 125.25
$ws.Rows.Item(6).RowHeight = 113.25

# Fill in the new "Goals Completed:" entries for each day.
# (set in this order so the shared-string table comes out in the same
# sequence as the authored workbook)
$ws.Range("D4").Value = "Did first weekly meeting form for Mike. Continued to do some reading on MVC 5"
$ws.Range("D5").Value = "Setup webservices environment`nTools in my hands now"
$ws.Range("D6").Value = "Keep testing the environment.`nGet key responses back from within my app.`nGet a logical response from webservices"
$ws.Range("D7").Value = "Rembering all the processes"
$ws.Range("D3").Value = "Gained access to my own WebServices environment`nStarted planning / structuring of the project`nAm now linked into their repository system so can start branching my project"

# Leave the selection where the author left it when saving.
$ws.Range("E3").Select()
